$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row: "_old" columns become "_FV2210", "_new" columns become "_FV2304".
$oldHeaders = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    # Left block: columns A..J (1..10) -> "_old" => "_FV2210"
    $ws.Cells.Item(1, $i + 1).Value = "$($oldHeaders[$i])_FV2210"
    # Right block: columns L..U (12..21) -> "_new" => "_FV2304"
    $ws.Cells.Item(1, $i + 12).Value = "$($oldHeaders[$i])_FV2304"
}
# Column K (11) stays "diff" - unchanged.

# 2. Turn the data range into a real Excel Table ("Table1") with an autofilter.
$tableRange = $ws.Range("A1:U61")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# 3. Freeze the header row (split/freeze below row 1).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
